# Belgium First Division B - base update (29-03-2024 17:05)
# Applies:
#   1) Ten pairs of adjacent rows whose full records (every column except
#      the running index in column A) were swapped.
#   2) A shift of rows 212-217 (cols B:AC) up into rows 210-215 (col A
#      stays put, i.e. the two records that used to occupy rows 210/211
#      were dropped), followed by deletion of the now-superfluous trailing
#      rows 216/217.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value = $v2
    $rng2.Value = $v1
}

# 1) Row-pair swaps (full record swap, column A/index untouched)
Swap-Rows 33 34
Swap-Rows 68 69
Swap-Rows 73 74
Swap-Rows 81 82
Swap-Rows 108 109
Swap-Rows 112 113
Swap-Rows 117 118
Swap-Rows 133 134
Swap-Rows 164 165
Swap-Rows 181 182

# 2) Drop the two records previously at rows 210 & 211: shift the B:AC
#    data of rows 212..217 up into rows 210..215 (column A keeps its
#    original sequential value in each destination row).
for ($i = 0; $i -lt 6; $i++) {
    $dst = 210 + $i
    $src = 212 + $i
    $srcRng = $ws.Range("B$src`:AC$src")
    $v = $srcRng.Value2
    $dstRng = $ws.Range("B$dst`:AC$dst")
    $dstRng.Value = $v
}

# The data that used to live in rows 216/217 has now been copied into
# rows 214/215, so the trailing rows are redundant - remove them so the
# sheet shrinks back down to 215 rows (dimension A1:AC215).
$ws.Range("A216:A217").EntireRow.Delete()
